$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently follows
#    the title ("Play Ghost Castle Free - Review of Cristaltec's Slot
#    Game 2021").
# ------------------------------------------------------------------
$metaFind = $d.Content
$metaFound = $metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($metaFound) {
    $metaPara = $metaFind.Paragraphs(1)
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Ghost Castle Free - Review of
#    Cristaltec's Slot Game 2021") right before the closing
#    image-prompt paragraph (the final, italicised paragraph of the
#    document).
# ------------------------------------------------------------------
$oldText = "Create a feature image that captures the essence of the game " + [char]34 + "Ghost Castle" + [char]34 + ". The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be standing in front of the Ghost Castle, which is shrouded in an eerie atmosphere. The warrior should be holding a stake with a hammer, ready to defeat any vampires that may come their way. The background should include a full moon, barren trees, and spooky ghosts floating around the castle. The overall feel of the image should be fun and engaging, while still capturing the haunting atmosphere of the game."
$newText = "Check out our review of Ghost Castle by Cristaltec, an online slot game with an eerie setting, spooky symbols, and a variety of game modes. Play free and win big!"

$imgFind = $d.Content
$imgFound = $imgFind.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$imgPara = $imgFind.Paragraphs(1)
$prevPara = $imgPara.Previous()

# Insert one character inside the previous paragraph's text (away from
# any paragraph boundary) so the runtime treats this as a plain
# insertion rather than a replace of a whole paragraph.
$insertAt = $prevPara.Range.End - 1
$insertRange = $d.Range($insertAt, $insertAt)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Ghost Castle Free - Review of Cristaltec''s Slot Game 2021</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 3) Swap the body text of the final (italic) paragraph from the old
#    image-generation prompt to the new meta-description text, while
#    keeping its italic formatting intact.
# ------------------------------------------------------------------
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
